# Rediseno red parte 2 - added vlans
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsVlans = $wb.Worksheets.Add($null, $lastSheet)
$wsVlans.Name = "VLANs"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsIps = $wb.Worksheets.Add($null, $lastSheet2)
$wsIps.Name = "Static IPs"

# ---------------------------------------------------------------------------
# 2. VLANs sheet content (table1 Tabla1 -> C5:H23)
# ---------------------------------------------------------------------------
$vlanData = @(
@('VLAN','Name','Network IP','Gateway','DHCP','HSRP'),
@(10,'WiFi_Empleados','192.168.10.0/24','192.168.10.1','192.168.10.10 - 192.168.10.49','Router_Movistar'),
@(20,'Demo_Area','192.168.20.0/24','192.168.20.1','192.168.20.10 - 192.168.20.49','Router_Movistar'),
@(30,'Developers','192.168.30.0/24','192.168.30.1','192.168.30.10 - 192.168.30.49','Router_Movistar'),
@(40,'Programers','192.168.40.0/24','192.168.40.1','192.168.40.10 - 192.168.40.49','Router_Movistar'),
@(49,'Native-Trunk','N/A','N/A','N/A','Router_Movistar'),
@(50,'Systems-administrators','192.168.50.0/24','192.168.50.1','192.168.50.10 - 192.168.50.49','Router_Movistar'),
@(60,'Commercial','192.168.60.0/24','192.168.60.1','192.168.60.10 - 192.168.60.49','Router_Movistar'),
@(70,'Accounting','192.168.70.0/24','192.168.70.1','192.168.60.10 - 192.168.60.49','Router_Movistar'),
@(80,'Call-Center','192.168.80.0/24','192.168.80.1','192.168.70.10 - 192.168.70.49','Router_Movistar'),
@(90,'R-H','192.168.90.0/24','192.168.90.1','192.168.90.10 - 192.168.90.49','Router_Movistar'),
@(99,'Management','192.168.99.0/24','N/A','N/A','Router_Movistar'),
@(100,'Direction','192.168.100.0/24','192.168.100.1','192.168.100.10 - 192.168.100.49','Router_Movistar'),
@(110,'Testing','192.168.110.0/24','192.168.110.1','N/A','Router_Movistar'),
@(120,'Printers','192.168.120.0/24','N/A','192.168.120.10 - 192.168.120.49','Router_Movistar'),
@(130,'Reception','192.168.130.0/24','192.168.130.1','192.168.130.10 - 192.168.130.15','Router_Movistar'),
@(150,'Servers','192.168.150.0/24','192.168.150.1','N/A','Router_Movistar'),
@(200,'DMZ','192.168.200.0/24','192.168.200.1','N/A','Router_Movistar'),
@(555,'Black-Hole','N/A','N/A','N/A',$null)
)

$startRow = 5
for ($i = 0; $i -lt $vlanData.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $vlanData[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $col = 3 + $c
        $v = $rowVals[$c]
        if ($null -ne $v) {
            $wsVlans.Cells.Item($row, $col).Value = $v
        }
    }
}

# Data cell formatting (rows 6-23, columns D:H) -> plain size 14, style index 9
$wsVlans.Range("D6:H23").Font.Size = 14

# Row heights
$wsVlans.Rows("5:23").RowHeight = 18.75

# Column widths (approximate "best fit" widths from the source file)
$wsVlans.Columns("C").ColumnWidth = 12.140625
$wsVlans.Columns("D").ColumnWidth = 27.85546875
$wsVlans.Columns("E").ColumnWidth = 21.85546875
$wsVlans.Columns("F").ColumnWidth = 17.85546875
$wsVlans.Columns("G").ColumnWidth = 39.7109375
$wsVlans.Columns("H").ColumnWidth = 20.85546875

# Table (Tabla1)
$tabla1 = $wsVlans.ListObjects.Add(1, $wsVlans.Range("C5:H23"), $null, 1)
$tabla1.Name = "Tabla1"
$tabla1.TableStyle = "TableStyleMedium4"

# ---------------------------------------------------------------------------
# 3. Static IPs sheet content (table2 Tabla2 -> D7:E15, title D6:E6)
# ---------------------------------------------------------------------------
$wsIps.Range("D6").Value = "VLAN 49 - Management"
$wsIps.Range("D6:E6").Merge() | Out-Null

$ipData = @(
@('Device Name','IP Address'),
@('Router_Movistar','192.168.49.50'),
@('Router_Vodafone','192.168.49.51'),
@('Switch_C1','192.168.49.100'),
@('Switch_C2','192.168.49.101'),
@('Switch_P1-1','192.168.49.102'),
@('Switch_P1-2','192.168.49.103'),
@('Switch_P1-3','192.168.49.104'),
@('Switch_p0-1','192.168.49.105')
)

$startRow2 = 7
for ($i = 0; $i -lt $ipData.Count; $i++) {
    $row = $startRow2 + $i
    $rowVals = $ipData[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $col = 4 + $c
        $wsIps.Cells.Item($row, $col).Value = $rowVals[$c]
    }
}

# Data formatting: size 14 -> style index 9
$wsIps.Range("D7:E15").Font.Size = 14

# Row heights
$wsIps.Rows("6").RowHeight = 28.5
$wsIps.Rows("7:15").RowHeight = 18.75

# Title formatting: size 22, centered + wrap -> style index 10
$wsIps.Range("D6:E6").Font.Size = 22
$wsIps.Range("D6:E6").HorizontalAlignment = -4108
$wsIps.Range("D6:E6").WrapText = $true

# VLANs header row (row 5) formatting: size 14, centered -> style index 11
$wsVlans.Range("C5:H5").Font.Size = 14
$wsVlans.Range("C5:H5").HorizontalAlignment = -4108

# VLAN id column (C6:C23) -> bold size 14, style index 12
$wsVlans.Range("C6:C23").Font.Size = 14
$wsVlans.Range("C6:C23").Font.Bold = $true

# Column widths
$wsIps.Columns("D").ColumnWidth = 20.85546875
$wsIps.Columns("E").ColumnWidth = 19.28515625

# Table (Tabla2)
$tabla2 = $wsIps.ListObjects.Add(1, $wsIps.Range("D7:E15"), $null, 1)
$tabla2.Name = "Tabla2"
$tabla2.TableStyle = "TableStyleMedium4"

# ---------------------------------------------------------------------------
# 4. Selections / active sheet
# ---------------------------------------------------------------------------
$wsVlans.Range("I22").Select()
$wsIps.Range("E26").Select()
$wsVlans.Activate()

Write-Host "done"
